$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# --- Row 3 (J3:P3) - latest Essential raw poll numbers (new poll added) ---
$ws.Range("J3").Value = 30
$ws.Range("K3").Value = 47
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 3
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 8
$ws.Range("P3").Value = 2
# Q3 is a formula (K3+L3*0.822+M3*0.348+N3*0.352+(O3+P3)*(0.507)) and will recalc automatically

# --- Rows 6-8 (B:G) - Essential poll swing table shifts down, new poll entered at row 6 ---
$ws.Range("B6").Value = 57
$ws.Range("C6").Value = 55
$ws.Range("D6").Value = 60.5
$ws.Range("E6").Value = 50.5
$ws.Range("F6").Value = 59
$ws.Range("G6").Value = 56

$ws.Range("B7").Value = 55.5
$ws.Range("C7").Value = 53
$ws.Range("D7").Value = 60
$ws.Range("E7").Value = 49
$ws.Range("F7").Value = 57
$ws.Range("G7").Value = 63.5

$ws.Range("B8").Value = 58
$ws.Range("C8").Value = 57.5
$ws.Range("D8").Value = 64
$ws.Range("E8").Value = 45.5
$ws.Range("F8").Value = 59
$ws.Range("G8").Value = 60.5

# --- Rows 9-11 (B:E, H) - Essential (3 avg.) rows shift down, new 3-avg entered at row 9 ---
$ws.Range("B9").Value = 54.878
$ws.Range("C9").Value = 52.726999999999997
$ws.Range("D9").Value = 53.233999999999995
$ws.Range("E9").Value = 53.542999999999999
$ws.Range("H9").Value = 60.86399999999999

$ws.Range("B10").Value = 53.277000000000001
$ws.Range("C10").Value = 53.266999999999996
$ws.Range("D10").Value = 53.224000000000004
$ws.Range("E10").Value = 52.902000000000001
$ws.Range("H10").Value = 54.402999999999992

$ws.Range("B11").Value = 53.704999999999998
$ws.Range("C11").Value = 52.332000000000001
$ws.Range("D11").Value = 58.515999999999998
$ws.Range("E11").Value = 48.938000000000002
$ws.Range("H11").Value = 54.055999999999997

# --- Rows 12-13 (B:G) - updated 3-average figures ---
$ws.Range("B12").Value = 52.499999999999993
$ws.Range("C12").Value = 51.236749116607776
$ws.Range("D12").Value = 51.27272727272728
$ws.Range("E12").Value = 52.158273381294968
$ws.Range("F12").Value = 53.191489361702125
$ws.Range("G12").Value = 52.631578947368418

$ws.Range("B13").Value = 52.142857142857139
$ws.Range("C13").Value = 49.635036496350367
$ws.Range("D13").Value = 52.857142857142861
$ws.Range("E13").Value = 51.957295373665481
$ws.Range("F13").Value = 54.255319148936174
$ws.Range("G13").Value = 50

# Recalculate all formulas (Q3, Q6, and the swing/deviation tables in rows 16-27)
$excel.CalculateFullRebuild()

# --- Update the active selection to match the saved view state ---
$ws.Activate()
$ws.Range("I20:I21").Select()
